$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grades for the first two tasks (columns C and D) for the rows
# that were still blank, plus a few corrections to already-graded rows.

$ws.Range("C2:D3").Value = 0

$ws.Range("C16:D18").Value = 0
$ws.Range("C19:D20").Value = 5
$ws.Range("C21:D22").Value = 0

$ws.Range("C24").Value = 5

$ws.Range("C25:D25").Value = 0

$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 0

$ws.Range("C27:D27").Value = 0
$ws.Range("C28:D28").Value = 5
$ws.Range("C29:D30").Value = 0

$ws.Range("C31").Value = 0

$ws.Range("D32").Value = 0
$ws.Range("C32").Value = 5
$ws.Range("C32").Font.Underline = 2
$ws.Range("C32").HorizontalAlignment = -4108

# Set the new column widths for C and D introduced alongside the grades.
$ws.Columns.Item(3).ColumnWidth = 7.333333333333333
$ws.Columns.Item(4).ColumnWidth = 5.5

# Move the active selection to where the editor left off.
[void]$ws.Range("D16").Select()
